$d = $word.ActiveDocument

# --- 1) Collapse the split "<id>" / "p012r_N" / "</id>" runs into a single run
#     "<id>p012r_N</id>" for each of the four <id> occurrences. Word's Find/Replace,
#     when the search hit spans multiple runs, merges the hit into one run that keeps
#     the formatting of the first run in the span (the brown Courier <id> formatting) -
#     exactly the formatting we want the merged run to use.

$ids = @("p012r_1", "p012r_2", "p012r_3", "p012r_4")
foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 1) | Out-Null
}

# --- 2) Wrap the <m>...</m> markup around "terre chimolée" with <def>...</def>,
#     i.e. "<m>" -> "<m><def>" and "</m>" -> "</def></m>", without merging those
#     runs with the "terre chimolée" run in between. We scope each Find to a Range
#     that starts right after the previous match so we hit only the one occurrence
#     (out of many "<m>"/"</m>" pairs in the document) that immediately follows
#     "La " and precedes "terre chimolée".

$total = $d.Content.End

$laRng = $d.Range(0, $total)
$laRng.Find.Execute("La ", $true) | Out-Null

$openRng = $d.Range($laRng.End, $total)
$openRng.Find.Execute("<m>", $true, $false, $false, $false, $false, $true, 1, $false, "<m><def>", 1) | Out-Null

$closeRng = $d.Range($openRng.End, $total)
$closeRng.Find.Execute("</m>", $true, $false, $false, $false, $false, $true, 1, $false, "</def></m>", 1) | Out-Null
